$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24; existing rows 24-107 shift down to 25-108.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new data record.
$ws.Cells.Item(24, 1).Value = 5
$ws.Cells.Item(24, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(24, 3).Value = "Maule"
$ws.Cells.Item(24, 4).Value = 45145
$ws.Cells.Item(24, 5).Value = 7
$ws.Cells.Item(24, 6).Value = 100112040
$ws.Cells.Item(24, 7).Value = "Cilantro"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 200
$ws.Cells.Item(24, 11).Value = 8000
$ws.Cells.Item(24, 12).Value = 8000
$ws.Cells.Item(24, 13).Value = 8000
$ws.Cells.Item(24, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(24, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 16).Value = 222
$ws.Cells.Item(24, 17).Value = 36
$ws.Cells.Item(24, 18).Value = "Hortaliza"
